$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price observation was recorded (weekly update). Insert a new data
# row right above the current row 202, shifting all subsequent rows (and
# the old last row 257) down by one — matching the diff, which shows every
# row from 202..257 taking on the values previously held by the row above
# it, with a brand-new row appearing at 202 and the old row 257's data
# ending up at the new row 258.
$ws.Rows.Item(202).Insert()

$ws.Range("A202").Value2 = 4
$ws.Range("B202").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C202").Value2 = "Los Lagos"
$ws.Range("D202").Value2 = 44932
$ws.Range("E202").Value2 = 10
$ws.Range("F202").Value2 = 100112009
$ws.Range("G202").Value2 = "Acelga"
$ws.Range("H202").Value2 = "Sin especificar"
$ws.Range("I202").Value2 = "Primera"
$ws.Range("J202").Value2 = 80
$ws.Range("K202").Value2 = 10000
$ws.Range("L202").Value2 = 10000
$ws.Range("M202").Value2 = 10000
$ws.Range("N202").Value2 = "`$/docena de atados (12 kilos)"
$ws.Range("O202").Value2 = "Región de La Araucanía"
$ws.Range("P202").Value2 = 833
$ws.Range("Q202").Value2 = 12
$ws.Range("R202").Value2 = "Hortaliza"
